# Realestate Update resale numbers 2023-06-23 09:14
# Appends one new data row (row 70) to the CityResaleNum sheet with the
# latest resale-number snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A and D hold values that look like dates/numbers ("2023-06-23",
# "25") but must stay plain text, matching the rest of the column. Force
# text via NumberFormat "@" before assigning, then ClearFormats so the
# cell is left with no explicit style (same as its neighbours in the
# original data).
$ws.Range("A70").NumberFormat = "@"
$ws.Range("A70").Value = "2023-06-23"
$ws.Range("A70").ClearFormats()

$ws.Range("B70").Value = "09:13:46"

$ws.Range("C70").Value = "Friday"

$ws.Range("D70").NumberFormat = "@"
$ws.Range("D70").Value = "25"
$ws.Range("D70").ClearFormats()

$ws.Range("E70").Value = 122462
$ws.Range("F70").Value = 133904
$ws.Range("G70").Value = 162281
$ws.Range("H70").Value = 133595
$ws.Range("I70").Value = 177334
$ws.Range("J70").Value = 115340
$ws.Range("K70").Value = 202312
$ws.Range("L70").Value = 225684
$ws.Range("M70").Value = 175527
$ws.Range("N70").Value = 104081
$ws.Range("O70").Value = 39361
$ws.Range("P70").Value = 33869
$ws.Range("Q70").Value = 51913
$ws.Range("R70").Value = -1
$ws.Range("S70").Value = 35622
$ws.Range("T70").Value = -1
